# Update F-column values ("view/heat" counts) on both the "展览" and
# "全部类型" worksheets to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 455
$ws1.Range("F3").Value = 5493
$ws1.Range("F5").Value = 63
$ws1.Range("F6").Value = 81
$ws1.Range("F10").Value = 13

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 455
$ws4.Range("F3").Value = 5493
$ws4.Range("F6").Value = 63
$ws4.Range("F7").Value = 81
$ws4.Range("F12").Value = 13
